# Weekly update: insert a new price record (row 176) for
# Terminal Hortofrutícola Agro Chillán - Mango, pushing the
# existing rows 176-202 down to 177-203.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 176 (shifts 176..202 -> 177..203)
$ws.Rows.Item(176).Insert()

# Populate the new row 176 with the new observation
$ws.Cells.Item(176, 1).Value = 7
$ws.Cells.Item(176, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(176, 3).Value = 'Ñuble'
$ws.Cells.Item(176, 4).Value = 45218
$ws.Cells.Item(176, 5).Value = 16
$ws.Cells.Item(176, 6).Value = 'Fruta'
$ws.Cells.Item(176, 7).Value = 100108
$ws.Cells.Item(176, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(176, 9).Value = 100108002
$ws.Cells.Item(176, 10).Value = 'Mango'
$ws.Cells.Item(176, 11).Value = 'Sin especificar'
$ws.Cells.Item(176, 12).Value = 'Primera'
$ws.Cells.Item(176, 13).Value = 30
$ws.Cells.Item(176, 14).Value = 10000
$ws.Cells.Item(176, 15).Value = 10000
$ws.Cells.Item(176, 16).Value = 10000
$ws.Cells.Item(176, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(176, 18).Value = 'Brasil'
$ws.Cells.Item(176, 19).Value = 2500
$ws.Cells.Item(176, 20).Value = 4
